$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are numeric-looking text (e.g. "42.802.60", "0.0900") that
# must stay as literal text, so force the cell format to Text before writing.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.802.60"
$ws.Range("E2").Value = "  -0.70%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.237.77"
$ws.Range("E3").Value = "  -1.65%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "114.44"
$ws.Range("E5").Value = "  +2.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "283.52"
$ws.Range("E6").Value = "  +7.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.622"
$ws.Range("E7").Value = "  -4.03%  "
$ws.Range("E8").Value = "  -0.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.612"
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.60"
$ws.Range("E10").Value = "  +0.30%  "
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.14"
$ws.Range("E12").Value = "  -1.26%  "
$ws.Range("E13").Value = "  -2.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.38"
$ws.Range("E14").Value = "  +0.91%  "
$ws.Range("E15").Value = "  +2.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.573.82"
$ws.Range("E16").Value = "  -1.66%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.242.35"
$ws.Range("E17").Value = "  -1.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.631.17"
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000107"
$ws.Range("E19").Value = "  -0.70%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.88"
$ws.Range("E20").Value = "  +2.44%  "
$ws.Range("E21").Value = "  +0.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.16"
$ws.Range("E22").Value = "  +10.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.35"
$ws.Range("E23").Value = "  -2.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "231.37"
$ws.Range("E24").Value = "  -1.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.17"
$ws.Range("E25").Value = "  -1.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.02"
$ws.Range("E26").Value = "  +5.97%  "
$ws.Range("E27").Value = "  -1.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.94"
$ws.Range("E28").Value = "  -0.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.12"
$ws.Range("E29").Value = "  -1.75%  "
$ws.Range("E30").Value = "  -1.57%  "
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "175.61"
$ws.Range("E32").Value = "  +1.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.14"
$ws.Range("E33").Value = "  -1.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0900"
$ws.Range("E34").Value = "  +0.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.63"
$ws.Range("E35").Value = "  +18.89%  "
$ws.Range("E36").Value = "  -0.90%  "
$ws.Range("E37").Value = "  -2.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0372"
$ws.Range("E38").Value = "  -1.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.64"
$ws.Range("E39").Value = "  -0.46%  "
$ws.Range("E40").Value = "  +1.48%  "
$ws.Range("E41").Value = "  +1.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.40"
$ws.Range("E42").Value = "  -2.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.42"
$ws.Range("E43").Value = "  -5.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.233"
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("E46").Value = "  -1.65%  "
$ws.Range("E47").Value = "  -7.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.29"
$ws.Range("E48").Value = "  +2.31%  "
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.651"
$ws.Range("E50").Value = "  +8.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "101.17"
$ws.Range("E51").Value = "  +1.56%  "
